# Update handback report timestamps for the "defb06e1..." row (row 3)
# in both the zh-cn and de-de worksheets, reflecting freshly generated
# handoff/handback datetimes.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-26 09:39:21"
$wsZhCn.Range("G3").Value = "2016-01-26 09:40:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-26 09:39:33"
$wsDeDe.Range("G3").Value = "2016-01-26 09:40:26"
